$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.207.69"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.578.08"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.97%  "
$ws.Range("D5").Value = "213.28"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "23.54"
$ws.Range("E8").Value = "  +6.56%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "1.802.88"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.573.28"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "28.157.57"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "63.82"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "230.44"
$ws.Range("E18").Value = "  +6.42%  "
$ws.Range("D19").Value = "0.0₃0707"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "7.47"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "152.21"
$ws.Range("D26").Value = "15.26"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "6.60"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").Value = "1.416.39"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "2.52"
$ws.Range("E39").Value = "  +7.43%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "0.972"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").Value = "1.82"
$ws.Range("E45").Value = "  +4.67%  "
$ws.Range("D46").Value = "63.94"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "1.715.12"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "87.13"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("E51").Value = "  -1.65%  "
